$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 6: IsDeleted (F6) changes from 0 to 1
$ws.Range("F6").Value = 1

# Add new row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 19
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = "2025-03-27 17:47:33"
$ws.Range("E7").Value = 8000.08
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = "O1047"

# Add new row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 18
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = "2025-03-27 17:50:33"
$ws.Range("E8").Value = 15000.15
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = "O1050"
